$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 4 new rows before row 95, shifting existing rows 95:140 down to 99:144
$ws.Range("A95:A98").EntireRow.Insert()

# Row 95: Lane Late / Primera
$ws.Cells.Item(95, 1).Value = 11
$ws.Cells.Item(95, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(95, 3).Value = "Bíobío"
$ws.Cells.Item(95, 4).Value = 44468
$ws.Cells.Item(95, 5).Value = 8
$ws.Cells.Item(95, 6).Value = "Fruta"
$ws.Cells.Item(95, 7).Value = 100102
$ws.Cells.Item(95, 8).Value = "Cítricos"
$ws.Cells.Item(95, 9).Value = 100102005
$ws.Cells.Item(95, 10).Value = "Naranja"
$ws.Cells.Item(95, 11).Value = "Lane Late"
$ws.Cells.Item(95, 12).Value = "Primera"
$ws.Cells.Item(95, 13).Value = 100
$ws.Cells.Item(95, 14).Value = 7000
$ws.Cells.Item(95, 15).Value = 7500
$ws.Cells.Item(95, 16).Value = 7250
$ws.Cells.Item(95, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(95, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(95, 19).Value = 483
$ws.Cells.Item(95, 20).Value = 15

# Row 96: Lane Late / Segunda
$ws.Cells.Item(96, 1).Value = 11
$ws.Cells.Item(96, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(96, 3).Value = "Bíobío"
$ws.Cells.Item(96, 4).Value = 44468
$ws.Cells.Item(96, 5).Value = 8
$ws.Cells.Item(96, 6).Value = "Fruta"
$ws.Cells.Item(96, 7).Value = 100102
$ws.Cells.Item(96, 8).Value = "Cítricos"
$ws.Cells.Item(96, 9).Value = 100102005
$ws.Cells.Item(96, 10).Value = "Naranja"
$ws.Cells.Item(96, 11).Value = "Lane Late"
$ws.Cells.Item(96, 12).Value = "Segunda"
$ws.Cells.Item(96, 13).Value = 50
$ws.Cells.Item(96, 14).Value = 6500
$ws.Cells.Item(96, 15).Value = 6500
$ws.Cells.Item(96, 16).Value = 6500
$ws.Cells.Item(96, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(96, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(96, 19).Value = 433
$ws.Cells.Item(96, 20).Value = 15

# Row 97: Navel Late / Primera
$ws.Cells.Item(97, 1).Value = 11
$ws.Cells.Item(97, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(97, 3).Value = "Bíobío"
$ws.Cells.Item(97, 4).Value = 44468
$ws.Cells.Item(97, 5).Value = 8
$ws.Cells.Item(97, 6).Value = "Fruta"
$ws.Cells.Item(97, 7).Value = 100102
$ws.Cells.Item(97, 8).Value = "Cítricos"
$ws.Cells.Item(97, 9).Value = 100102005
$ws.Cells.Item(97, 10).Value = "Naranja"
$ws.Cells.Item(97, 11).Value = "Navel Late"
$ws.Cells.Item(97, 12).Value = "Primera"
$ws.Cells.Item(97, 13).Value = 100
$ws.Cells.Item(97, 14).Value = 7000
$ws.Cells.Item(97, 15).Value = 7500
$ws.Cells.Item(97, 16).Value = 7250
$ws.Cells.Item(97, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(97, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(97, 19).Value = 483
$ws.Cells.Item(97, 20).Value = 15

# Row 98: Navel Late / Segunda
$ws.Cells.Item(98, 1).Value = 11
$ws.Cells.Item(98, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(98, 3).Value = "Bíobío"
$ws.Cells.Item(98, 4).Value = 44468
$ws.Cells.Item(98, 5).Value = 8
$ws.Cells.Item(98, 6).Value = "Fruta"
$ws.Cells.Item(98, 7).Value = 100102
$ws.Cells.Item(98, 8).Value = "Cítricos"
$ws.Cells.Item(98, 9).Value = 100102005
$ws.Cells.Item(98, 10).Value = "Naranja"
$ws.Cells.Item(98, 11).Value = "Navel Late"
$ws.Cells.Item(98, 12).Value = "Segunda"
$ws.Cells.Item(98, 13).Value = 50
$ws.Cells.Item(98, 14).Value = 6500
$ws.Cells.Item(98, 15).Value = 6500
$ws.Cells.Item(98, 16).Value = 6500
$ws.Cells.Item(98, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(98, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(98, 19).Value = 433
$ws.Cells.Item(98, 20).Value = 15
